$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("ZZ1").Value = 1
$ws.Range("ZZ1").Style = "Normal"
$ws.Range("ZZ2").Value = 2
$ws.Range("ZZ3").Value = 3
$ws.Range("ZZ3").Style = "Comma"
Write-Host "done"
